$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.960542333333334
$ws.Range("H2").Value = 14.881627
$ws.Range("I2").Value = 0.1435881646191863
$ws.Range("J2").Value = 0.1435881646191863
$ws.Range("M2").Value = 6.139446
$ws.Range("N2").Value = 18.418338
$ws.Range("O2").Value = 0.1067674894682656
$ws.Range("P2").Value = 0.1067674894682656
$ws.Range("Q2").Value = 30.454981786214
$ws.Range("R2").Value = 274.094836075926
$ws.Range("S2").Value = 0.01533054785374656
$ws.Range("T2").Value = 0.01533054785374656
$ws.Range("G3").Value = 4.960542333333334
$ws.Range("H3").Value = 14.881627
$ws.Range("I3").Value = 0.1435881646191863
$ws.Range("J3").Value = 0.1435881646191863
$ws.Range("O3").Value = 0.3167683832774248
$ws.Range("P3").Value = 0.3167683832774247
$ws.Range("Q3").Value = 90.35686229214789
$ws.Range("R3").Value = 813.2117606293311
$ws.Range("S3").Value = 0.04548419076419237
$ws.Range("T3").Value = 0.04548419076419236
$ws.Range("G4").Value = 4.960542333333334
$ws.Range("H4").Value = 14.881627
$ws.Range("I4").Value = 0.1435881646191863
$ws.Range("J4").Value = 0.1435881646191863
$ws.Range("M4").Value = 13.317528
$ws.Range("N4").Value = 39.95258399999999
$ws.Range("O4").Value = 0.231597285892462
$ws.Range("P4").Value = 0.231597285892462
$ws.Range("Q4").Value = 66.06216141935199
$ws.Range("R4").Value = 594.559452774168
$ws.Range("S4").Value = 0.03325462921208359
$ws.Range("T4").Value = 0.03325462921208359
$ws.Range("G5").Value = 4.960542333333334
$ws.Range("H5").Value = 14.881627
$ws.Range("I5").Value = 0.1435881646191863
$ws.Range("J5").Value = 0.1435881646191863
$ws.Range("M5").Value = 13.58033733333333
$ws.Range("N5").Value = 40.741012
$ws.Range("O5").Value = 0.2361676482230093
$ws.Range("P5").Value = 0.2361676482230092
$ws.Range("Q5").Value = 67.36583824294711
$ws.Range("R5").Value = 606.2925441865241
$ws.Range("S5").Value = 0.03391087915077153
$ws.Range("T5").Value = 0.03391087915077153
$ws.Range("G6").Value = 4.960542333333334
$ws.Range("H6").Value = 14.881627
$ws.Range("I6").Value = 0.1435881646191863
$ws.Range("J6").Value = 0.1435881646191863
$ws.Range("M6").Value = 6.250524666666666
$ws.Range("N6").Value = 18.751574
$ws.Range("O6").Value = 0.1086991931388382
$ws.Range("P6").Value = 0.1086991931388382
$ws.Range("Q6").Value = 31.00599221454422
$ws.Range("R6").Value = 279.053929930898
$ws.Range("S6").Value = 0.01560791763839223
$ws.Range("T6").Value = 0.01560791763839223
$ws.Range("I7").Value = 0.1393077946862016
$ws.Range("J7").Value = 0.1393077946862016
$ws.Range("M7").Value = 6.139446
$ws.Range("N7").Value = 18.418338
$ws.Range("O7").Value = 0.1067674894682656
$ws.Range("P7").Value = 0.1067674894682656
$ws.Range("Q7").Value = 29.54711734841
$ws.Range("R7").Value = 265.92405613569
$ws.Range("S7").Value = 0.01487354350200634
$ws.Range("T7").Value = 0.01487354350200634
$ws.Range("I8").Value = 0.1393077946862016
$ws.Range("J8").Value = 0.1393077946862016
$ws.Range("O8").Value = 0.3167683832774248
$ws.Range("P8").Value = 0.3167683832774247
$ws.Range("S8").Value = 0.04412830490069153
$ws.Range("T8").Value = 0.04412830490069151
$ws.Range("I9").Value = 0.1393077946862016
$ws.Range("J9").Value = 0.1393077946862016
$ws.Range("M9").Value = 13.317528
$ws.Range("N9").Value = 39.95258399999999
$ws.Range("O9").Value = 0.231597285892462
$ws.Range("P9").Value = 0.231597285892462
$ws.Range("Q9").Value = 64.09284528388
$ws.Range("R9").Value = 576.83560755492
$ws.Range("S9").Value = 0.03226330715298864
$ws.Range("T9").Value = 0.03226330715298864
$ws.Range("I10").Value = 0.1393077946862016
$ws.Range("J10").Value = 0.1393077946862016
$ws.Range("M10").Value = 13.58033733333333
$ws.Range("N10").Value = 40.741012
$ws.Range("O10").Value = 0.2361676482230093
$ws.Range("P10").Value = 0.2361676482230092
$ws.Range("Q10").Value = 65.35765944011779
$ws.Range("R10").Value = 588.21893496106
$ws.Range("S10").Value = 0.03289999425017406
$ws.Range("T10").Value = 0.03289999425017406
$ws.Range("I11").Value = 0.1393077946862016
$ws.Range("J11").Value = 0.1393077946862016
$ws.Range("M11").Value = 6.250524666666666
$ws.Range("N11").Value = 18.751574
$ws.Range("O11").Value = 0.1086991931388382
$ws.Range("P11").Value = 0.1086991931388382
$ws.Range("Q11").Value = 30.08170212998555
$ws.Range("R11").Value = 270.73531916987
$ws.Range("S11").Value = 0.01514264488034106
$ws.Range("T11").Value = 0.01514264488034105
$ws.Range("G12").Value = 11.78248366666667
$ws.Range("H12").Value = 35.347451
$ws.Range("I12").Value = 0.3410564996056291
$ws.Range("J12").Value = 0.3410564996056292
$ws.Range("M12").Value = 6.139446
$ws.Range("N12").Value = 18.418338
$ws.Range("O12").Value = 0.1067674894682656
$ws.Range("P12").Value = 0.1067674894682656
$ws.Range("Q12").Value = 72.33792221738199
$ws.Range("R12").Value = 651.0412999564379
$ws.Range("S12").Value = 0.03641374622972755
$ws.Range("T12").Value = 0.03641374622972755
$ws.Range("G13").Value = 11.78248366666667
$ws.Range("H13").Value = 35.347451
$ws.Range("I13").Value = 0.3410564996056291
$ws.Range("J13").Value = 0.3410564996056292
$ws.Range("O13").Value = 0.3167683832774248
$ws.Range("P13").Value = 0.3167683832774247
$ws.Range("Q13").Value = 214.6193263939115
$ws.Range("R13").Value = 1931.573937545203
$ws.Range("S13").Value = 0.1080359159863328
$ws.Range("T13").Value = 0.1080359159863328
$ws.Range("G14").Value = 11.78248366666667
$ws.Range("H14").Value = 35.347451
$ws.Range("I14").Value = 0.3410564996056291
$ws.Range("J14").Value = 0.3410564996056292
$ws.Range("M14").Value = 13.317528
$ws.Range("N14").Value = 39.95258399999999
$ws.Range("O14").Value = 0.231597285892462
$ws.Range("P14").Value = 0.231597285892462
$ws.Range("Q14").Value = 156.913556140376
$ws.Range("R14").Value = 1412.222005263384
$ws.Range("S14").Value = 0.07898775964464726
$ws.Range("T14").Value = 0.07898775964464727
$ws.Range("G15").Value = 11.78248366666667
$ws.Range("H15").Value = 35.347451
$ws.Range("I15").Value = 0.3410564996056291
$ws.Range("J15").Value = 0.3410564996056292
$ws.Range("M15").Value = 13.58033733333333
$ws.Range("N15").Value = 40.741012
$ws.Range("O15").Value = 0.2361676482230093
$ws.Range("P15").Value = 0.2361676482230092
$ws.Range("Q15").Value = 160.0101028178235
$ws.Range("R15").Value = 1440.090925360412
$ws.Range("S15").Value = 0.08054651142303311
$ws.Range("T15").Value = 0.08054651142303312
$ws.Range("G16").Value = 11.78248366666667
$ws.Range("H16").Value = 35.347451
$ws.Range("I16").Value = 0.3410564996056291
$ws.Range("J16").Value = 0.3410564996056292
$ws.Range("M16").Value = 6.250524666666666
$ws.Range("N16").Value = 18.751574
$ws.Range("O16").Value = 0.1086991931388382
$ws.Range("P16").Value = 0.1086991931388382
$ws.Range("Q16").Value = 73.6467047930971
$ws.Range("R16").Value = 662.8203431378739
$ws.Range("S16").Value = 0.03707256632188839
$ws.Range("T16").Value = 0.03707256632188839
$ws.Range("G17").Value = 0.9139316666666666
$ws.Range("H17").Value = 2.741795
$ws.Range("I17").Value = 0.02645472244480135
$ws.Range("J17").Value = 0.02645472244480136
$ws.Range("M17").Value = 6.139446
$ws.Range("N17").Value = 18.418338
$ws.Range("O17").Value = 0.1067674894682656
$ws.Range("P17").Value = 0.1067674894682656
$ws.Range("Q17").Value = 5.611034115189999
$ws.Range("R17").Value = 50.49930703670999
$ws.Range("S17").Value = 0.002824504300011218
$ws.Range("T17").Value = 0.002824504300011219
$ws.Range("G18").Value = 0.9139316666666666
$ws.Range("H18").Value = 2.741795
$ws.Range("I18").Value = 0.02645472244480135
$ws.Range("J18").Value = 0.02645472244480136
$ws.Range("O18").Value = 0.3167683832774248
$ws.Range("P18").Value = 0.3167683832774247
$ws.Range("Q18").Value = 16.64737284762611
$ws.Range("R18").Value = 149.826355628635
$ws.Range("S18").Value = 0.008380019658892728
$ws.Range("T18").Value = 0.008380019658892726
$ws.Range("G19").Value = 0.9139316666666666
$ws.Range("H19").Value = 2.741795
$ws.Range("I19").Value = 0.02645472244480135
$ws.Range("J19").Value = 0.02645472244480136
$ws.Range("M19").Value = 13.317528
$ws.Range("N19").Value = 39.95258399999999
$ws.Range("O19").Value = 0.231597285892462
$ws.Range("P19").Value = 0.231597285892462
$ws.Range("Q19").Value = 12.17131056092
$ws.Range("R19").Value = 109.54179504828
$ws.Range("S19").Value = 0.00612684191725439
$ws.Range("T19").Value = 0.006126841917254391
$ws.Range("G20").Value = 0.9139316666666666
$ws.Range("H20").Value = 2.741795
$ws.Range("I20").Value = 0.02645472244480135
$ws.Range("J20").Value = 0.02645472244480136
$ws.Range("M20").Value = 13.58033733333333
$ws.Range("N20").Value = 40.741012
$ws.Range("O20").Value = 0.2361676482230093
$ws.Range("P20").Value = 0.2361676482230092
$ws.Range("Q20").Value = 12.41150033294889
$ws.Range("R20").Value = 111.70350299654
$ws.Range("S20").Value = 0.006247749584181192
$ws.Range("T20").Value = 0.006247749584181192
$ws.Range("G21").Value = 0.9139316666666666
$ws.Range("H21").Value = 2.741795
$ws.Range("I21").Value = 0.02645472244480135
$ws.Range("J21").Value = 0.02645472244480136
$ws.Range("M21").Value = 6.250524666666666
$ws.Range("N21").Value = 18.751574
$ws.Range("O21").Value = 0.1086991931388382
$ws.Range("P21").Value = 0.1086991931388382
$ws.Range("Q21").Value = 5.712552426147777
$ws.Range("R21").Value = 51.41297183532999
$ws.Range("S21").Value = 0.002875606984461821
$ws.Range("T21").Value = 0.002875606984461821
$ws.Range("G22").Value = 12.077388
$ws.Range("H22").Value = 36.232164
$ws.Range("I22").Value = 0.3495928186441815
$ws.Range("J22").Value = 0.3495928186441815
$ws.Range("M22").Value = 6.139446
$ws.Range("N22").Value = 18.418338
$ws.Range("O22").Value = 0.1067674894682656
$ws.Range("P22").Value = 0.1067674894682656
$ws.Range("Q22").Value = 74.148471447048
$ws.Range("R22").Value = 667.3362430234318
$ws.Range("S22").Value = 0.03732514758277394
$ws.Range("T22").Value = 0.03732514758277394
$ws.Range("G23").Value = 12.077388
$ws.Range("H23").Value = 36.232164
$ws.Range("I23").Value = 0.3495928186441815
$ws.Range("J23").Value = 0.3495928186441815
$ws.Range("O23").Value = 0.3167683832774248
$ws.Range("P23").Value = 0.3167683832774247
$ws.Range("Q23").Value = 219.991043525988
$ws.Range("R23").Value = 1979.919391733892
$ws.Range("S23").Value = 0.1107399519673153
$ws.Range("T23").Value = 0.1107399519673153
$ws.Range("G24").Value = 12.077388
$ws.Range("H24").Value = 36.232164
$ws.Range("I24").Value = 0.3495928186441815
$ws.Range("J24").Value = 0.3495928186441815
$ws.Range("M24").Value = 13.317528
$ws.Range("N24").Value = 39.95258399999999
$ws.Range("O24").Value = 0.231597285892462
$ws.Range("P24").Value = 0.231597285892462
$ws.Range("Q24").Value = 160.840952856864
$ws.Range("R24").Value = 1447.568575711776
$ws.Range("S24").Value = 0.08096474796548812
$ws.Range("T24").Value = 0.08096474796548812
$ws.Range("G25").Value = 12.077388
$ws.Range("H25").Value = 36.232164
$ws.Range("I25").Value = 0.3495928186441815
$ws.Range("J25").Value = 0.3495928186441815
$ws.Range("M25").Value = 13.58033733333333
$ws.Range("N25").Value = 40.741012
$ws.Range("O25").Value = 0.2361676482230093
$ws.Range("P25").Value = 0.2361676482230092
$ws.Range("Q25").Value = 164.015003145552
$ws.Range("R25").Value = 1476.135028309968
$ws.Range("S25").Value = 0.08256251381484932
$ws.Range("T25").Value = 0.08256251381484932
$ws.Range("G26").Value = 12.077388
$ws.Range("H26").Value = 36.232164
$ws.Range("I26").Value = 0.3495928186441815
$ws.Range("J26").Value = 0.3495928186441815
$ws.Range("M26").Value = 6.250524666666666
$ws.Range("N26").Value = 18.751574
$ws.Range("O26").Value = 0.1086991931388382
$ws.Range("P26").Value = 0.1086991931388382
$ws.Range("Q26").Value = 75.49001160290399
$ws.Range("R26").Value = 679.4101044261358
$ws.Range("S26").Value = 0.03800045731375473
$ws.Range("T26").Value = 0.03800045731375473
